$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 428; this shifts rows 428..522 down to 429..523
# and carries formatting (e.g. the date style on column D) along with it.
$ws.Rows.Item(428).Insert()

# Populate the newly inserted row 428 with the new data point.
# Columns that are identical to the (now shifted-down) former row 428 are
# simply re-entered with the same values; D/J/N/O/P/Q carry the new values.
$ws.Range("A428").Value = 5
$ws.Range("B428").Value = "Macroferia Regional de Talca"
$ws.Range("C428").Value = "Maule"
$ws.Range("D428").Value = 44889
$ws.Range("E428").Value = 7
$ws.Range("F428").Value = 100112043
$ws.Range("G428").Value = "Pepino ensalada"
$ws.Range("H428").Value = "Sin especificar"
$ws.Range("I428").Value = "Primera"
$ws.Range("J428").Value = 500
$ws.Range("K428").Value = 16000
$ws.Range("L428").Value = 16000
$ws.Range("M428").Value = 16000
$ws.Range("N428").Value = "$/caja 80 unidades"
$ws.Range("O428").Value = "Región del Maule"
$ws.Range("P428").Value = 200
$ws.Range("Q428").Value = 80
$ws.Range("R428").Value = "Hortaliza"
